$wb = $excel.ActiveWorkbook

# ALC!row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 162.35715
$ws.Range("I33").Value = 187.1
$ws.Range("J33").Value = 100.5
$ws.Range("K33").Value = 187.1
$ws.Range("L33").Value = 100.5
$ws.Range("M33").Value = 41.90000000000001
$ws.Range("N33").Value = -558.5

# ALC!row 103
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H103").Value = 7692895
$ws.Range("I103").Value = 420
$ws.Range("J103").Value = 15385370
$ws.Range("K103").Value = 1260
$ws.Range("L103").Value = 46156110
$ws.Range("M103").Value = -674
$ws.Range("N103").Value = -46157282

# ALC!row 111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 2880.0625
$ws.Range("I111").Value = 2474.4443
$ws.Range("J111").Value = 3401.5715
$ws.Range("K111").Value = 7423.3329
$ws.Range("L111").Value = 10204.7145
$ws.Range("M111").Value = -4356.3329
$ws.Range("N111").Value = -16338.7145

# ALC!row 124
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H124").Value = 20657.857
$ws.Range("J124").Value = 20657.857
$ws.Range("L124").Value = 20657.857
$ws.Range("N124").Value = -30477.857

# ALC!row 125
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 1488.4445
$ws.Range("I125").Value = 483.33334
$ws.Range("J125").Value = 1689.4667
$ws.Range("K125").Value = 4350.00006
$ws.Range("L125").Value = 15205.2003
$ws.Range("M125").Value = -1890.00006
$ws.Range("N125").Value = -20125.2003

# ALC!row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 35867.777
$ws.Range("J130").Value = 35867.777
$ws.Range("L130").Value = 35867.777
$ws.Range("N130").Value = -45907.777

# ALC!row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 346587.22
$ws.Range("I135").Value = 346587.22
$ws.Range("K135").Value = 3119284.98
$ws.Range("M135").Value = -3116749.98

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2140.3076
$ws.Range("I137").Value = 1266.8572
$ws.Range("J137").Value = 4363.636
$ws.Range("K137").Value = 3800.5716
$ws.Range("L137").Value = 13090.908
$ws.Range("M137").Value = -1250.5716
$ws.Range("N137").Value = -18190.908

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2457.2207
$ws.Range("I138").Value = 2429.3044
$ws.Range("J138").Value = 2471.4888
$ws.Range("K138").Value = 7287.9132
$ws.Range("L138").Value = 7414.4664
$ws.Range("M138").Value = -2147.9132
$ws.Range("N138").Value = -17694.4664

# ARM!row 3
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 4321.25
$ws.Range("J3").Value = 4560
$ws.Range("L3").Value = 4560
$ws.Range("N3").Value = -4790

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1889.8125
$ws.Range("I61").Value = 1870.9231
$ws.Range("J61").Value = 1971.6666
$ws.Range("K61").Value = 1870.9231
$ws.Range("L61").Value = 1971.6666
$ws.Range("M61").Value = -1658.9231
$ws.Range("N61").Value = -2395.6666

# ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1384.6428
$ws.Range("I74").Value = 721.3111
$ws.Range("J74").Value = 4098.273
$ws.Range("K74").Value = 721.3111
$ws.Range("L74").Value = 4098.273
$ws.Range("M74").Value = 152.6889
$ws.Range("N74").Value = -5846.273

# ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1384.6428
$ws.Range("I77").Value = 721.3111
$ws.Range("J77").Value = 4098.273
$ws.Range("K77").Value = 3606.5555
$ws.Range("L77").Value = 20491.365
$ws.Range("M77").Value = 761.4445000000001
$ws.Range("N77").Value = -29227.365

# ARM!row 81
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# ARM!row 84
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1889.8125
$ws.Range("I136").Value = 1870.9231
$ws.Range("J136").Value = 1971.6666
$ws.Range("K136").Value = 5612.7693
$ws.Range("L136").Value = 5914.9998
$ws.Range("M136").Value = -3062.7693
$ws.Range("N136").Value = -11014.9998

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2774.6667
$ws.Range("I31").Value = 1398.6428
$ws.Range("J31").Value = 4256.5386
$ws.Range("K31").Value = 1398.6428
$ws.Range("L31").Value = 4256.5386
$ws.Range("M31").Value = -1103.6428
$ws.Range("N31").Value = -4846.5386

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2774.6667
$ws.Range("I34").Value = 1398.6428
$ws.Range("J34").Value = 4256.5386
$ws.Range("K34").Value = 1398.6428
$ws.Range("L34").Value = 4256.5386
$ws.Range("M34").Value = -1196.6428
$ws.Range("N34").Value = -4660.5386

# CRP!row 81
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()

# CRP!row 84
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()

# CRP!row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1904.0834
$ws.Range("I99").Value = 1756.125
$ws.Range("J99").Value = 2200
$ws.Range("K99").Value = 1756.125
$ws.Range("L99").Value = 2200
$ws.Range("M99").Value = -258.125
$ws.Range("N99").Value = -5196

# CRP!row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1904.0834
$ws.Range("I126").Value = 1756.125
$ws.Range("J126").Value = 2200
$ws.Range("K126").Value = 5268.375
$ws.Range("L126").Value = 6600
$ws.Range("M126").Value = -2798.375
$ws.Range("N126").Value = -11540

# CRP!row 127
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 54588
$ws.Range("J127").Value = 54485
$ws.Range("L127").Value = 54485
$ws.Range("N127").Value = -64405

# CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 8883.362999999999
$ws.Range("I134").Value = 9786
$ws.Range("J134").Value = 3166.6667
$ws.Range("K134").Value = 29358
$ws.Range("L134").Value = 9500.000100000001
$ws.Range("M134").Value = -26823
$ws.Range("N134").Value = -14570.0001

# CUL!row 14
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 285.94736
$ws.Range("I14").Value = 285.94736
$ws.Range("K14").Value = 857.84208
$ws.Range("M14").Value = -684.84208

# GSM!row 48
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value = 1000
$ws.Range("I48").Value = 1000
$ws.Range("K48").Value = 1000
$ws.Range("M48").Value = -515

# GSM!row 74
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# GSM!row 77
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# GSM!row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2922.8572
$ws.Range("J80").Value = 3526.6667
$ws.Range("L80").Value = 3526.6667
$ws.Range("N80").Value = -5522.6667

# GSM!row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2922.8572
$ws.Range("J83").Value = 3526.6667
$ws.Range("L83").Value = 17633.3335
$ws.Range("N83").Value = -27617.3335

# GSM!row 123
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 34361.875
$ws.Range("J123").Value = 34361.875
$ws.Range("L123").Value = 34361.875
$ws.Range("N123").Value = -39261.875

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2359.2666
$ws.Range("I132").Value = 2011.28
$ws.Range("J132").Value = 4099.2
$ws.Range("K132").Value = 6033.84
$ws.Range("L132").Value = 12297.6
$ws.Range("M132").Value = -3503.84
$ws.Range("N132").Value = -17357.6

# GSM!row 135
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 53143.2
$ws.Range("J135").Value = 53143.2
$ws.Range("L135").Value = 53143.2
$ws.Range("N135").Value = -63283.2

# LTW!row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1997.2941
$ws.Range("I61").Value = 1304.1538
$ws.Range("J61").Value = 4250
$ws.Range("K61").Value = 1304.1538
$ws.Range("L61").Value = 4250
$ws.Range("M61").Value = -1102.1538
$ws.Range("N61").Value = -4654

# LTW!row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1997.2941
$ws.Range("I113").Value = 1304.1538
$ws.Range("J113").Value = 4250
$ws.Range("K113").Value = 1304.1538
$ws.Range("L113").Value = 4250
$ws.Range("M113").Value = 865.8462
$ws.Range("N113").Value = -8590

# LTW!row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2633.238
$ws.Range("I122").Value = 1706.125
$ws.Range("J122").Value = 5600
$ws.Range("K122").Value = 5118.375
$ws.Range("L122").Value = 16800
$ws.Range("M122").Value = -2668.375
$ws.Range("N122").Value = -21700

# WVR!row 54
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 7077
$ws.Range("J54").Value = 7077
$ws.Range("L54").Value = 7077
$ws.Range("N54").Value = -8117

# WVR!row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 7117.579
$ws.Range("I81").Value = 14172.875
$ws.Range("J81").Value = 1986.4546
$ws.Range("K81").Value = 28345.75
$ws.Range("L81").Value = 3972.9092
$ws.Range("M81").Value = -27284.75
$ws.Range("N81").Value = -6094.9092

# WVR!row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 7117.579
$ws.Range("I84").Value = 14172.875
$ws.Range("J84").Value = 1986.4546
$ws.Range("K84").Value = 141728.75
$ws.Range("L84").Value = 19864.546
$ws.Range("M84").Value = -136424.75
$ws.Range("N84").Value = -30472.546

# WVR!row 109
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 26218
$ws.Range("J109").Value = 26218
$ws.Range("L109").Value = 26218
$ws.Range("N109").Value = -28992

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1498.4584
$ws.Range("I136").Value = 1452.8636
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4358.5908
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1808.5908
$ws.Range("N136").Value = -11100
